{"js": "// Replace each multiplication expression with its new value.\n// The mapping below preserves document order (each old value is unique),\n// so a direct search-and-replace per pair is safe and unambiguous.\nconst replacements = [\n  [\"268\u00d75=\", \"270\u00d75=\"],\n  [\"794\u00d79=\", \"749\u00d74=\"],\n  [\"621\u00d77=\", \"824\u00d74=\"],\n  [\"250\u00d76=\", \"986\u00d77=\"],\n  [\"679\u00d72=\", \"883\u00d75=\"],\n  [\"624\u00d78=\", \"625\u00d72=\"],\n  [\"529\u00d79=\", \"596\u00d72=\"],\n  [\"466\u00d75=\", \"780\u00d73=\"],\n  [\"336\u00d78=\", \"890\u00d77=\"],\n  [\"536\u00d74=\", \"980\u00d73=\"],\n  [\"291\u00d73=\", \"938\u00d73=\"],\n  [\"613\u00d72=\", \"189\u00d79=\"],\n  [\"257\u00d75=\", \"315\u00d74=\"],\n  [\"593\u00d74=\", \"125\u00d77=\"],\n  [\"597\u00d74=\", \"431\u00d73=\"],\n  [\"275\u00d75=\", \"966\u00d75=\"],\n  [\"112\u00d76=\", \"255\u00d76=\"],\n  [\"856\u00d79=\", \"536\u00d75=\"],\n  [\"456\u00d76=\", \"939\u00d75=\"],\n  [\"439\u00d75=\", \"413\u00d77=\"],\n  [\"922\u00d72=\", \"897\u00d73=\"],\n  [\"354\u00d74=\", \"149\u00d74=\"],\n  [\"736\u00d76=\", \"912\u00d79=\"],\n  [\"582\u00d72=\", \"889\u00d72=\"],\n  [\"150\u00d74=\", \"647\u00d72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each multiplication expression to its new value.\n# Each \"old\" string is unique in the document, so Find/Replace per pair\n# is unambiguous and preserves all run formatting (font, size, etc.).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"268\u00d75=\"; New = \"270\u00d75=\" }\n    @{ Old = \"794\u00d79=\"; New = \"749\u00d74=\" }\n    @{ Old = \"621\u00d77=\"; New = \"824\u00d74=\" }\n    @{ Old = \"250\u00d76=\"; New = \"986\u00d77=\" }\n    @{ Old = \"679\u00d72=\"; New = \"883\u00d75=\" }\n    @{ Old = \"624\u00d78=\"; New = \"625\u00d72=\" }\n    @{ Old = \"529\u00d79=\"; New = \"596\u00d72=\" }\n    @{ Old = \"466\u00d75=\"; New = \"780\u00d73=\" }\n    @{ Old = \"336\u00d78=\"; New = \"890\u00d77=\" }\n    @{ Old = \"536\u00d74=\"; New = \"980\u00d73=\" }\n    @{ Old = \"291\u00d73=\"; New = \"938\u00d73=\" }\n    @{ Old = \"613\u00d72=\"; New = \"189\u00d79=\" }\n    @{ Old = \"257\u00d75=\"; New = \"315\u00d74=\" }\n    @{ Old = \"593\u00d74=\"; New = \"125\u00d77=\" }\n    @{ Old = \"597\u00d74=\"; New = \"431\u00d73=\" }\n    @{ Old = \"275\u00d75=\"; New = \"966\u00d75=\" }\n    @{ Old = \"112\u00d76=\"; New = \"255\u00d76=\" }\n    @{ Old = \"856\u00d79=\"; New = \"536\u00d75=\" }\n    @{ Old = \"456\u00d76=\"; New = \"939\u00d75=\" }\n    @{ Old = \"439\u00d75=\"; New = \"413\u00d77=\" }\n    @{ Old = \"922\u00d72=\"; New = \"897\u00d73=\" }\n    @{ Old = \"354\u00d74=\"; New = \"149\u00d74=\" }\n    @{ Old = \"736\u00d76=\"; New = \"912\u00d79=\" }\n    @{ Old = \"582\u00d72=\"; New = \"889\u00d72=\" }\n    @{ Old = \"150\u00d74=\"; New = \"647\u00d72=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair.New\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
